$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.262.20"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "3.025.43"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'541.36"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'133.43"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.023.24"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.492"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "'6.16"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D11").Value = "'0.147"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "'34.25"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "3.510.61"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "62.190.40"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "3.021.73"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("D19").Value = "'6.63"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").Value = "'480.20"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("D21").Value = "'13.24"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'0.672"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").Value = "'80.63"
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("D25").Value = "'12.08"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'2.70"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("D28").Value = "'7.70"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("D31").Value = "'25.68"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "'1.13"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").Value = "'5.64"
$ws.Range("E33").Value = "  +5.09%  "
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("D35").Value = "'54.98"
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("D36").Value = "'5.86"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "'452.76"
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D38").Value = "3.161.03"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").Value = "'8.09"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "'2.46"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").Value = "'26.33"
$ws.Range("E44").Value = "  +5.72%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "'1.96"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "'114.12"
$ws.Range("E49").Value = "  -5.99%  "
$ws.Range("D50").Value = "0.0₃0496"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("E51").Value = "  +4.38%  "
